$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 242, shifting existing rows 242-267 down to 243-268.
$ws.Rows.Item(242).Insert()

# Populate the newly inserted row 242 with the new weekly record.
$ws.Cells.Item(242, 1).Value = 3
$ws.Cells.Item(242, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(242, 3).Value = "Coquimbo"
$ws.Cells.Item(242, 4).Value = 44449
$ws.Cells.Item(242, 5).Value = 5
$ws.Cells.Item(242, 6).Value = 100112037
$ws.Cells.Item(242, 7).Value = "Cebollín"
$ws.Cells.Item(242, 8).Value = "Sin especificar"
$ws.Cells.Item(242, 9).Value = "Primera"
$ws.Cells.Item(242, 10).Value = 300
$ws.Cells.Item(242, 11).Value = 3500
$ws.Cells.Item(242, 12).Value = 3700
$ws.Cells.Item(242, 13).Value = 3607
$ws.Cells.Item(242, 14).Value = "`$/paquete 36 unidades"
$ws.Cells.Item(242, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(242, 16).Value = 100
$ws.Cells.Item(242, 17).Value = 36
$ws.Cells.Item(242, 18).Value = "Hortaliza"
